$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 'M. Cenk Gursoy'
$ws.Cells.Item(2, 8).Value = 'Svetoslava Todorova'
$ws.Cells.Item(3, 7).Value = 'Jay Henderson'
$ws.Cells.Item(3, 8).Value = 'Yaoying Wu'
$ws.Cells.Item(4, 7).Value = 'Jay Henderson'
$ws.Cells.Item(4, 8).Value = 'Bing Dong'
$ws.Cells.Item(5, 7).Value = 'Chikukuri Monhan'
$ws.Cells.Item(5, 8).Value = 'Anupam Pandey'
$ws.Cells.Item(6, 7).Value = 'Elizabeth Carter'
$ws.Cells.Item(6, 8).Value = 'Sucheta Soundarajan'
$ws.Cells.Item(7, 7).Value = 'Yiyang Sun'
$ws.Cells.Item(7, 8).Value = 'Yi Zheng'
$ws.Cells.Item(8, 7).Value = 'Anupam Pandey'
$ws.Cells.Item(8, 8).Value = 'Zhenyu Gan'
$ws.Cells.Item(9, 7).Value = 'Yaoying Wu'
$ws.Cells.Item(9, 8).Value = 'Zhen Ma'
$ws.Cells.Item(10, 7).Value = 'Bing Dong'
$ws.Cells.Item(10, 8).Value = 'Svetoslava Todorova'
$ws.Cells.Item(11, 7).Value = 'Svetoslava Todorova'
$ws.Cells.Item(11, 8).Value = 'Shalabh Maroo'
$ws.Cells.Item(12, 7).Value = 'Zhenyu Gan'
$ws.Cells.Item(12, 8).Value = 'Shalabh Maroo'
$ws.Cells.Item(13, 7).Value = 'Farzana Rahman'
$ws.Cells.Item(13, 8).Value = 'C.Y. Roger Chen'
$ws.Cells.Item(14, 7).Value = 'Shalabh Maroo'
$ws.Cells.Item(14, 8).Value = 'Zhenyu Gan'
$ws.Cells.Item(15, 7).Value = 'Zhen Ma'
$ws.Cells.Item(15, 8).Value = 'Ashok Sangani'
$ws.Cells.Item(16, 7).Value = 'Ruth Chen'
$ws.Cells.Item(16, 8).Value = 'Chikukuri Monhan'
$ws.Cells.Item(17, 7).Value = 'Bing Dong'
$ws.Cells.Item(17, 8).Value = 'Baris Salman'
$ws.Cells.Item(18, 7).Value = 'Zhenyu Gan'
$ws.Cells.Item(18, 8).Value = 'Sucheta Soundarajan'
$ws.Cells.Item(19, 7).Value = 'Yaoying Wu'
$ws.Cells.Item(19, 8).Value = 'Gabriel Silva De Oliveira'
$ws.Cells.Item(20, 7).Value = 'Jason Pollack'
$ws.Cells.Item(20, 8).Value = 'Nadeem Ghani'
$ws.Cells.Item(21, 7).Value = 'Ben Akih-Kumgeh'
$ws.Cells.Item(21, 8).Value = 'Ben Akih-Kumgeh'
$ws.Cells.Item(22, 7).Value = 'Amit Sanyal'
$ws.Cells.Item(22, 8).Value = 'Elizabeth Carter'
$ws.Cells.Item(23, 7).Value = 'Ruth Chen'
$ws.Cells.Item(23, 8).Value = 'Ruth Chen'
$ws.Cells.Item(24, 7).Value = 'Mary Beth Monroe'
$ws.Cells.Item(24, 8).Value = 'Bing Dong'
$ws.Cells.Item(25, 7).Value = 'Jason Pollack'
$ws.Cells.Item(25, 8).Value = 'Yi Zheng'
$ws.Cells.Item(26, 7).Value = 'Ashok Sangani'
$ws.Cells.Item(26, 8).Value = 'Amit Sanyal'
$ws.Cells.Item(27, 7).Value = 'Jesse Q. Bond'
$ws.Cells.Item(27, 8).Value = 'Ashok Sangani'
$ws.Cells.Item(28, 7).Value = 'Younes Radi'
$ws.Cells.Item(28, 8).Value = 'Elizabeth Carter'
$ws.Cells.Item(29, 7).Value = 'Zhenyu Gan'
$ws.Cells.Item(29, 8).Value = 'M. Cenk Gursoy'
$ws.Cells.Item(30, 7).Value = 'John F. Dannenhoffer'
$ws.Cells.Item(30, 8).Value = 'Yuzhe Tang'
$ws.Cells.Item(31, 7).Value = 'Yiyang Sun'
$ws.Cells.Item(31, 8).Value = 'Wanliang Shan'
$ws.Cells.Item(32, 7).Value = 'Svetoslava Todorova'
$ws.Cells.Item(32, 8).Value = 'Senem Velipasalar'
$ws.Cells.Item(33, 7).Value = 'Pankaj Jha'
$ws.Cells.Item(33, 8).Value = 'Shikha Nangia'
$ws.Cells.Item(34, 7).Value = 'M. Cenk Gursoy'
$ws.Cells.Item(34, 8).Value = 'Gabriel Silva De Oliveira'
$ws.Cells.Item(35, 7).Value = 'Shikha Nangia'
$ws.Cells.Item(35, 8).Value = 'C.Y. Roger Chen'
$ws.Cells.Item(36, 7).Value = 'Ben Akih-Kumgeh'
$ws.Cells.Item(36, 8).Value = 'Anupam Pandey'
$ws.Cells.Item(37, 7).Value = 'Anupam Pandey'
$ws.Cells.Item(37, 8).Value = 'John F. Dannenhoffer'
$ws.Cells.Item(38, 7).Value = 'Chikukuri Monhan'
$ws.Cells.Item(38, 8).Value = 'Joao Paulo Marum'
$ws.Cells.Item(39, 7).Value = 'Elizabeth Carter'
$ws.Cells.Item(39, 8).Value = 'Endadul Hoque'
$ws.Cells.Item(40, 7).Value = 'C.Y. Roger Chen'
$ws.Cells.Item(40, 8).Value = 'Qinru Qiu'
$ws.Cells.Item(41, 7).Value = 'Gabriel Silva De Oliveira'
$ws.Cells.Item(41, 8).Value = 'Ruth Chen'
$ws.Cells.Item(42, 7).Value = 'Mary Beth Monroe'
$ws.Cells.Item(42, 8).Value = 'Shikha Nangia'
$ws.Cells.Item(43, 7).Value = 'Shikha Nangia'
$ws.Cells.Item(43, 8).Value = 'Pankaj Jha'
$ws.Cells.Item(44, 7).Value = 'Mary Beth Monroe'
$ws.Cells.Item(44, 8).Value = 'Jeongmin Ahn'
$ws.Cells.Item(45, 7).Value = 'Min Liu'
$ws.Cells.Item(45, 8).Value = 'Yi Zheng'
$ws.Cells.Item(46, 7).Value = 'Yi Zheng'
$ws.Cells.Item(46, 8).Value = 'Nadeem Ghani'
$ws.Cells.Item(47, 7).Value = 'Yaoying Wu'
$ws.Cells.Item(47, 8).Value = 'Jay Henderson'
$ws.Cells.Item(48, 7).Value = 'Anupam Pandey'
$ws.Cells.Item(48, 8).Value = 'Baris Salman'
$ws.Cells.Item(49, 7).Value = 'Jeongmin Ahn'
$ws.Cells.Item(49, 8).Value = 'Bing Dong'
$ws.Cells.Item(50, 7).Value = 'Elizabeth Carter'
$ws.Cells.Item(50, 8).Value = 'Chikukuri Monhan'
$ws.Cells.Item(51, 7).Value = 'Senem Velipasalar'
$ws.Cells.Item(51, 8).Value = 'Elizabeth Carter'
$ws.Cells.Item(52, 7).Value = 'Ashok Sangani'
$ws.Cells.Item(52, 8).Value = 'Min Liu'
$ws.Cells.Item(53, 7).Value = 'Era Jain'
$ws.Cells.Item(53, 8).Value = 'Baris Salman'
$ws.Cells.Item(54, 7).Value = 'Joao Paulo Marum'
$ws.Cells.Item(54, 8).Value = 'Shikha Nangia'
$ws.Cells.Item(55, 7).Value = 'Chikukuri Monhan'
$ws.Cells.Item(55, 8).Value = 'Qinru Qiu'
$ws.Cells.Item(56, 7).Value = 'Yuzhe Tang'
$ws.Cells.Item(56, 8).Value = 'Younes Radi'
$ws.Cells.Item(57, 7).Value = 'Ian Hosein'
$ws.Cells.Item(57, 8).Value = 'Baris Salman'
$ws.Cells.Item(58, 7).Value = 'Jesse Q. Bond'
$ws.Cells.Item(58, 8).Value = 'Ian Hosein'
$ws.Cells.Item(59, 7).Value = 'Baris Salman'
$ws.Cells.Item(59, 8).Value = 'Jesse Q. Bond'
$ws.Cells.Item(60, 7).Value = 'Mary Beth Monroe'
$ws.Cells.Item(60, 8).Value = 'Joao Paulo Marum'
$ws.Cells.Item(61, 7).Value = 'Ashok Sangani'
$ws.Cells.Item(61, 8).Value = 'Ruth Chen'
$ws.Cells.Item(62, 7).Value = 'Zhenyu Gan'
$ws.Cells.Item(62, 8).Value = 'Sucheta Soundarajan'
$ws.Cells.Item(63, 7).Value = 'Bing Dong'
$ws.Cells.Item(63, 8).Value = 'Mary Beth Monroe'
$ws.Cells.Item(64, 7).Value = 'Endadul Hoque'
$ws.Cells.Item(64, 8).Value = 'Senem Velipasalar'
$ws.Cells.Item(65, 7).Value = 'Ruth Chen'
$ws.Cells.Item(65, 8).Value = 'Era Jain'
$ws.Cells.Item(66, 7).Value = 'John F. Dannenhoffer'
$ws.Cells.Item(66, 8).Value = 'Pankaj Jha'
$ws.Cells.Item(67, 7).Value = 'Nadeem Ghani'
$ws.Cells.Item(67, 8).Value = 'Farzana Rahman'
$ws.Cells.Item(68, 7).Value = 'Jeongmin Ahn'
$ws.Cells.Item(68, 8).Value = 'Yiyang Sun'
$ws.Cells.Item(69, 7).Value = 'Wanliang Shan'
$ws.Cells.Item(69, 8).Value = 'Senem Velipasalar'
